$d = $word.ActiveDocument
$nl = [char]11

function Set-ParaPlainText($Index, $Text) {
    $p = $d.Paragraphs.Item($Index)
    $start = $p.Range.Start
    $end = $p.Range.End - 1
    $rng = $d.Range($start, $end)
    $rng.Text = $Text
}

# --- Paragraph 6 (Objetivos body): gets the old "Programa resumido" body text ---
$text6 = "Cinemática do corpo rígido" + $nl + "Dinâmica do ponto" + $nl + "Dinâmica do corpo rígido " + $nl + "Introdução às vibrações mecânicas"
Set-ParaPlainText 6 $text6

# --- Paragraph 8 (Docente(s) body, ListBullet): gets the old "Objetivos" body text ---
$text8 = "Proporcionar ao aluno conhecimento básico e compreensão de cinemática e dinâmica do corpo rígido. Desenvolver algumas aplicações práticas com ênfase em problemas bidimensionais. Apresentar conceitos fundamentais e exemplos das vibrações mecânicas."
Set-ParaPlainText 8 $text8

# --- Paragraph 10 (Programa resumido body): gets the old "Programa" body text ---
$text10 = "Cinemática do corpo rígido:" + $nl + `
    "Aceleração e velocidade angulares. Vínculo e cinemática do corpo rígido. Rotação em torno de um eixo fixo. Movimento plano e centro de rotação. Composição de movimentos. Composição de movimentos de rotação." + $nl + `
    "Dinâmica do ponto:" + $nl + `
    "Princípios da dinâmica do ponto. Teorema da resultante. Teorema da energia cinética para partícula. Teorema da quantidade de movimento." + $nl + `
    "Dinâmica do corpo rígido:" + $nl + `
    "Teorema do movimento do baricentro. Teorema da energia cinética para um sistema de partículas. Teorema do momento angular para um sistema de partículas. Teorema da energia cinética para o corpo rígido. Teorema do momento angular para corpo rígido Exercícios de aplicação: problemas bidimensionais. Rotação do corpo rígido, Balanceamento. Movimento de um giroscópio." + $nl + `
    "Introdução às vibrações mecânicas:" + $nl + `
    "Vibrações de sistemas mecânicos com um grau de liberdade: livres sem amortecimento, livres com amortecimento, forçadas. Vibrações de sistemas mecânicos com dois e mais graus de liberdade. Exemplos."
Set-ParaPlainText 10 $text10

# --- Paragraph 12 (Programa body): becomes the short "A avaliação..." sentence ---
$text12 = "A avaliação será composta por duas provas (P1 e P2)."
Set-ParaPlainText 12 $text12

# --- Paragraph 14 (Avaliação body, ListBullet): reordered Método/Critério/Norma blocks
#     plus the bibliography list appended at the end ---
$metodoLabel = "Método: "
$criterioBlock = "NS = NP1+NP2; " + $nl + "NP1: questões da P1 valendo até 4p. no total; " + $nl + "NP2: questões da P2 valendo até 6 p. no total." + $nl
$criterioLabel = "Critério: "
$normaBlock = "A recuperação consistirá de uma prova de Recuperação (R), que irá compor a nota final (NF) da seguinte forma: NF = (R + NS)/2." + $nl
$normaLabel = "Norma de recuperação: "
$bibliography = "HIBBELER, R.C. Dinâmica - Mecânica para Engenharia. São Paulo: Pearson Brasil, 2011, 12ª ed., 608p. ISBN: 8576058146." + $nl + $nl + `
    "BEER, F.P., JOHNSTON Jr., E.R., CLAUSEN, W. E., Mecânica Vetorial para Engenheiros - Dinâmica, 7ª Edição, McGraw-Hill, São Paulo, 2006, 1355 p. " + $nl + $nl + `
    "FRANÇA, L. N. F., MATSUMURA, A. Z. Mecânica Geral. Edgard Blücher, 2001, 235 p." + $nl + $nl + `
    "SOTELO JR., J., FRANÇA, L.N.F., Introdução às vibrações mecânicas, Edgard Blücher, 2006, 168 p. ISBN: 9788521203384." + $nl + $nl + `
    "GREENWOOD, D. T. Principles of Dynamics. New York: Prentice-Hall, 2nd ed, 1988, 552 p." + $nl + $nl + `
    "TENENBAUM, R. A. Dinâmica. Editora UFRJ, 1997, 756 p." + $nl + $nl + `
    "GIACAGLIA, G. E., Mecânica Geral, Editora Campus, Rio de Janeiro, 1982."

$text14 = $metodoLabel + $criterioBlock + $criterioLabel + $normaBlock + $normaLabel + $bibliography
Set-ParaPlainText 14 $text14

# Now apply Bold to the three labels within paragraph 14
$p14 = $d.Paragraphs.Item(14)
$base = $p14.Range.Start

$pos = $base
$rngMetodo = $d.Range($pos, $pos + $metodoLabel.Length)
$rngMetodo.Font.Bold = 1
$pos = $pos + $metodoLabel.Length

$pos = $pos + $criterioBlock.Length
$rngCriterio = $d.Range($pos, $pos + $criterioLabel.Length)
$rngCriterio.Font.Bold = 1
$pos = $pos + $criterioLabel.Length

$pos = $pos + $normaBlock.Length
$rngNorma = $d.Range($pos, $pos + $normaLabel.Length)
$rngNorma.Font.Bold = 1

# --- Paragraph 16 (Bibliografia body): gets the old "7797767 - Viktor Pastoukhov" text ---
$text16 = "7797767 - Viktor Pastoukhov"
Set-ParaPlainText 16 $text16
